$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16
$ws.Cells.Item($row, 1).Value = "2025-05-01T11:47:45.694Z"
$ws.Cells.Item($row, 2).Value = "UNDP"
$ws.Cells.Item($row, 3).Value = "C4"
$ws.Cells.Item($row, 4).Value = "الرحلة 3"
$ws.Cells.Item($row, 5).Value = "الصمود"
$ws.Cells.Item($row, 6).Value = "يامن "
$ws.Cells.Item($row, 7).Value = "'421123"
$ws.Cells.Item($row, 8).Value = "'"
